$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings / header row (row 21) ---
# Column A: new "#" header, right aligned
$ws.Range("A21").Value = "#"
$ws.Range("A21").HorizontalAlignment = -4152

# E21 header text changes from "Time(ms)" to "Forward (ms)"
$ws.Range("E21").Value = "Forward (ms)"

# New F21 / G21 headers, matching style of existing header cells (B21:E21 - bold, centered)
$ws.Range("F21").Value = "Entailment (ms)"
$ws.Range("G21").Value = "AskZ3 (ms)"
$ws.Range("F21:G21").HorizontalAlignment = -4108
$ws.Range("F21:G21").Font.Bold = $true

# --- Updated values in existing rows (22-25), new E/F/G values (style s=1: centered) ---
$ws.Range("C22").Value = 35
$ws.Range("E22:G22").HorizontalAlignment = -4108
$ws.Range("E22").Value = 44.31
$ws.Range("F22").Value = 935.51
$ws.Range("G22").Value = 988.67

$ws.Range("C23").Value = 22
$ws.Range("E23:G23").HorizontalAlignment = -4108
$ws.Range("E23").Value = 26.32
$ws.Range("F23").Value = 180.21
$ws.Range("G23").Value = 205.1

$ws.Range("C24").Value = 25
$ws.Range("E24:G24").HorizontalAlignment = -4108
$ws.Range("E24").Value = 32.7
$ws.Range("F24").Value = 236.66
$ws.Range("G24").Value = 245.58

$ws.Range("C25").Value = 39
$ws.Range("E25:G25").HorizontalAlignment = -4108
$ws.Range("E25").Value = 67.36
$ws.Range("F25").Value = 406.18
$ws.Range("G25").Value = 404.73

# Rows 26-27: new empty F/G cells (same style as D/E i.e. s=1 centered)
$ws.Range("F26:G27").HorizontalAlignment = -4108

# Row 28: updated + new values
$ws.Range("C28").Value = 154
$ws.Range("D28").Value = 39
$ws.Range("E28:G28").HorizontalAlignment = -4108
$ws.Range("E28").Value = 44.52
$ws.Range("F28").Value = 795.5
$ws.Range("G28").Value = 735.34

# --- New row 29: Totol summary row (style s=1: centered, matches rest of table) ---
$ws.Range("B29:G29").HorizontalAlignment = -4108
$ws.Range("B29").Value = "Totol"
$ws.Range("C29").Formula = "=SUM(C22:C28)"

# --- Column widths for E, F, G ---
$ws.Columns.Item(5).ColumnWidth = 21.5
$ws.Columns.Item(6).ColumnWidth = 17.5
$ws.Columns.Item(7).ColumnWidth = 11.5

# --- Selection moves to E34 ---
$ws.Range("E34").Select()
